$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "__TEMP__"
$title1.Text = "Slide 1"

$caption1 = $s1.Shapes.Item(3).TextFrame.TextRange
$caption1.Text = "__TEMP__"
$caption1.Text = "an image"

$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "__TEMP__"
$title2.Text = "Slide 2"

$caption2 = $s2.Shapes.Item(4).TextFrame.TextRange
$caption2.Text = "__TEMP__"
$caption2.Text = "an image"
